$wb = $excel.ActiveWorkbook

# --- Sheet "uart": add two note cells in column E ---
$ws1 = $wb.Worksheets.Item("uart")
$ws1.Range("E3").Value = "test reset"
$ws1.Range("E9").Value = "has some error"
$null = $ws1.Range("E9").Select()

# --- Sheet "uart_rx": add one note cell in column E ---
$ws2 = $wb.Worksheets.Item("uart_rx")
$ws2.Range("E9").Value = "may be error"
$null = $ws2.Range("E9").Select()

# --- Sheet "uart_tx": add two note cells in column E ---
$ws3 = $wb.Worksheets.Item("uart_tx")
$ws3.Range("E6").Value = "一二三四五"
$ws3.Range("E10").Value = "上山打老虎"
$null = $ws3.Range("E20").Select()

# Leave the first sheet ("uart") as the active/selected sheet, matching the
# tabSelected flag already present in the source workbook.
$null = $ws1.Select()
$null = $ws1.Range("E9").Select()
